$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ficha técnica")

# Remove the "DIMENSIÓN" / "Accesibilidad" row from the "Ficha técnica" sheet
$ws.Rows.Item(3).Delete()

# Append the two new rows at the end of the table
$ws.Range("A7").Value = "TIPOIND"
$ws.Range("B7").Value = "Resultados"
$ws.Range("A8").Value = "CITA"
$ws.Range("B8").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"
